# Edit: "Creato TC_Specification per alcuni TC + apportate modifiche al dizionario dei dati"
# The DD_Azi (Dizionario Dati Azienda) block in Foglio1 is corrected:
#   - C12 (formato dati for "nome"):            Stringa -> Valori alfanumerici
#   - row 16 "provincia": length rule corrected, and the description typo
#       "composta ad due caratteri" -> "composta da due caratteri"
#   - row 18 "partita iva": length rule and format corrected
#       "11 caratteri" -> "esattamente 11 caratteri"
#       Formato dati: Valori alfanumerici -> Numero

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

$ws.Range("C12").Value = "Valori alfanumerici"

$ws.Range("B16").Value = "esattamente 2 caratteri"
$ws.Range("D16").Value = "Provincia della città, composta da due caratteri"

$ws.Range("B18").Value = "esattamente 11 caratteri"
$ws.Range("C18").Value = "Numero"

# Restore the view state: scroll position and active cell selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B18").Select()
